# v0.6.2p: Can add/change terrain in GameMapEditor
# Adds new terrain entries (Sidewalk, Sand, Water) to the "Terrains" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terrains")

# --- Insert 7 new rows at row 17 for Sidewalk (132-133) and Sand (141-145) ---
$ws.Range("17:23").Insert()

$ws.Cells.Item(17, 2).Value = 132
$ws.Cells.Item(17, 3).Value = "Sidewalk, smooth"

$ws.Cells.Item(18, 2).Value = 133
$ws.Cells.Item(18, 3).Value = "Sidewalk, cracked"

$ws.Cells.Item(19, 2).Value = 141
$ws.Cells.Item(19, 3).Value = "Sand, light"

$ws.Cells.Item(20, 2).Value = 142
$ws.Cells.Item(20, 3).Value = "Sand, tan"

$ws.Cells.Item(21, 2).Value = 143
$ws.Cells.Item(21, 3).Value = "Sand, dark"

$ws.Cells.Item(22, 2).Value = 144
$ws.Cells.Item(22, 3).Value = "Sand, line left"

$ws.Cells.Item(23, 2).Value = 145
$ws.Cells.Item(23, 3).Value = "Sand, line up"

# --- Insert 5 new rows at row 42 for Water (181-185) ---
$ws.Range("42:46").Insert()

$ws.Cells.Item(42, 2).Value = 181
$ws.Cells.Item(42, 3).Value = "Water, rocks"

$ws.Cells.Item(43, 2).Value = 182
$ws.Cells.Item(43, 3).Value = "Water, dirt"

$ws.Cells.Item(44, 2).Value = 183
$ws.Cells.Item(44, 3).Value = "Water, shallow"

$ws.Cells.Item(45, 2).Value = 184
$ws.Cells.Item(45, 3).Value = "Water, medium"

$ws.Cells.Item(46, 2).Value = 185
$ws.Cells.Item(46, 3).Value = "Water, deep"

# --- Update the view to match where the edits were made ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E43").Select()
